$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived data for rows 2-5 (FAPs/MuSCs sending clusters x ECs/MuSCs targets).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric NATMI metrics.
$data = @(
    @("FAPs","Col2a1","Mag","ECs",1,0.3333333333333333,0.2021233333333333,0.60637,0.9764223557676824,0.9764223557676823,1,0.3333333333333333,0.00535,0.01605,0.003591913026022235,0.003591913026022235,0.001081359833333333,0.009732238499999999,0.003507224178581255,0.003507224178581255),
    @("FAPs","Col2a1","Mag","MuSCs",1,0.3333333333333333,0.2021233333333333,0.60637,0.9764223557676824,0.9764223557676823,3,1,1.484107,4.452321,0.9964080869739778,0.9964080869739778,0.2999726538633333,2.699753884769999,0.9729151315891011,0.972915131589101),
    @("MuSCs","Col2a1","Mag","ECs",2,0.6666666666666666,0.004880666666666667,0.014642,0.02357764423231757,0.02357764423231757,1,0.3333333333333333,0.00535,0.01605,0.003591913026022235,0.003591913026022235,0.00002611156666666667,0.0002350041,0.00008468884744097952,0.00008468884744097952),
    @("MuSCs","Col2a1","Mag","MuSCs",2,0.6666666666666666,0.004880666666666667,0.014642,0.02357764423231757,0.02357764423231757,3,1,1.484107,4.452321,0.9964080869739778,0.9964080869739778,0.007243431564666667,0.06519088408199999,0.02349295538487659,0.02349295538487659)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Old rows 6 and 7 are gone entirely now (used-range shrinks from A1:T7 to A1:T5).
$ws.Rows.Item(6).Resize(2).Delete() | Out-Null
